$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the formatting of the existing
# header cell H1 (bold font, borders, centered/top alignment) so the new
# headers look consistent with the rest of the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Add the new data values in row 2 (plain numeric cells, no special style)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
